# Fruta / hortaliza, semanal
# Insert a new weekly record as row 13 (Repollo, Primera calidad),
# shifting the previously existing rows 13-18 down to rows 14-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 13; this pushes rows
# 13-18 down to 14-19 and extends the sheet dimension to A1:R19.
$ws.Range("A13").EntireRow.Insert()

# Populate the newly inserted row 13 with the new data point.
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44799
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 100112006
$ws.Range("G13").Value = "Repollo"
$ws.Range("H13").Value = "Copenhague"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = 1100
$ws.Range("N13").Value = "`$/unidad"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 1100
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
